$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.207776069641113
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.202883839607239
